# Update "想去人数" (F column) figures across all sheets to match refreshed data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1204
$ws.Range("F6").Value = 9421
$ws.Range("F7").Value = 7265
$ws.Range("F8").Value = 196
$ws.Range("F10").Value = 5787
$ws.Range("F12").Value = 80
$ws.Range("F13").Value = 25
$ws.Range("F14").Value = 6657
$ws.Range("F15").Value = 1118
$ws.Range("F16").Value = 474
$ws.Range("F17").Value = 444
$ws.Range("F19").Value = 653
$ws.Range("F21").Value = 292
$ws.Range("F25").Value = 10789
$ws.Range("F26").Value = 95
$ws.Range("F27").Value = 38
$ws.Range("F28").Value = 2048
$ws.Range("F29").Value = 2573
$ws.Range("F32").Value = 2371
$ws.Range("F33").Value = 91
$ws.Range("F35").Value = 30
$ws.Range("F37").Value = 325
$ws.Range("F38").Value = 1489
$ws.Range("F40").Value = 22
$ws.Range("F41").Value = 5492
$ws.Range("F42").Value = 1224
$ws.Range("F43").Value = 770
$ws.Range("F44").Value = 139
$ws.Range("F45").Value = 173
$ws.Range("F47").Value = 1440
$ws.Range("F48").Value = 76
$ws.Range("F49").Value = 1111

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 102
$ws.Range("F20").Value = 31

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 103
$ws.Range("F3").Value = 185

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1204
$ws.Range("F5").Value = 9421
$ws.Range("F6").Value = 7265
$ws.Range("F7").Value = 103
$ws.Range("F8").Value = 196
$ws.Range("F9").Value = 185
$ws.Range("F12").Value = 5787
$ws.Range("F13").Value = 80
$ws.Range("F15").Value = 25
$ws.Range("F16").Value = 6657
$ws.Range("F17").Value = 6657
$ws.Range("F18").Value = 1118
$ws.Range("F19").Value = 474
$ws.Range("F20").Value = 444
$ws.Range("F21").Value = 653
$ws.Range("F23").Value = 292
$ws.Range("F28").Value = 10789
$ws.Range("F29").Value = 95
$ws.Range("F30").Value = 38
$ws.Range("F31").Value = 2048
$ws.Range("F32").Value = 2573
$ws.Range("F33").Value = 2371
$ws.Range("F34").Value = 91
$ws.Range("F36").Value = 30
$ws.Range("F38").Value = 325
$ws.Range("F39").Value = 1489
$ws.Range("F40").Value = 5492
$ws.Range("F41").Value = 31
$ws.Range("F42").Value = 1224
$ws.Range("F43").Value = 770
$ws.Range("F44").Value = 139
$ws.Range("F45").Value = 173
$ws.Range("F48").Value = 1440
$ws.Range("F49").Value = 76
$ws.Range("F50").Value = 1111

